$wb = $excel.ActiveWorkbook

# ---- Sheet: Summary ----
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.6086142322097379
$ws1.Range("C2").Value = 0.5650224215246636
$ws1.Range("D2").Value = 0.9438202247191011
$ws1.Range("E2").Value = 0.7068723702664796
$ws1.Range("F2").Value = 0.8322324966974901
$ws1.Range("G2").Value = 0.9200954922061508
$ws1.Range("H2").Value = 0.7761207900237064
$ws1.Range("I2").Value = 504
$ws1.Range("J2").Value = 388
$ws1.Range("K2").Value = 146
$ws1.Range("L2").Value = 30

# ---- Sheet: Classification Report ----
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.8295454545454546
$ws2.Range("C2").Value = 0.2734082397003745
$ws2.Range("D2").Value = 0.4112676056338028

$ws2.Range("B3").Value = 0.5650224215246636
$ws2.Range("C3").Value = 0.9438202247191011
$ws2.Range("D3").Value = 0.7068723702664796

$ws2.Range("B4").Value = 0.6086142322097379
$ws2.Range("C4").Value = 0.6086142322097379
$ws2.Range("D4").Value = 0.6086142322097379
$ws2.Range("E4").Value = 0.6086142322097379

$ws2.Range("B5").Value = 0.6972839380350591
$ws2.Range("C5").Value = 0.6086142322097379
$ws2.Range("D5").Value = 0.5590699879501413

$ws2.Range("B6").Value = 0.6972839380350591
$ws2.Range("C6").Value = 0.6086142322097379
$ws2.Range("D6").Value = 0.5590699879501412

# ---- Sheet: Confusion Matrix ----
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 146
$ws3.Range("C2").Value = 388
$ws3.Range("B3").Value = 30
$ws3.Range("C3").Value = 504
